$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (columns D, M, N, O, P, S) resulting from re-sequencing the
# weekly price records. Values taken from the corresponding source row
# before the edit.
$rows = @(
    @{ Row = 2; D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 },
    @{ Row = 3; D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  S = 629 },
    @{ Row = 4; D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 },
    @{ Row = 5; D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  S = 676 },
    @{ Row = 6; D = 44229; M = 55; N = 11000; O = 12000; P = 11364; S = 812 },
    @{ Row = 7; D = 44253; M = 90; N = 12000; O = 13000; P = 12667; S = 905 },
    @{ Row = 8; D = 44232; M = 60; N = 11000; O = 12000; P = 11583; S = 827 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D   # D - Fecha
    $ws.Cells.Item($r.Row, 13).Value = $r.M  # M - Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N  # N - Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O  # O - Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P  # P - Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S  # S - Precio $/Kg
}
